$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7170929312705994
$ws.Range("B1").Value = 1.366512656211853
$ws.Range("C1").Value = 4.258912563323975
$ws.Range("D1").Value = 2.036332845687866
$ws.Range("E1").Value = 0.8857404589653015
